$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the billing address row values (row 2)
$ws.Range("A2").Value = "Kollapudi"
$ws.Range("B2").Value = "Venu"
$ws.Range("C2").Value = "venukollapudi@gmail.com"
$ws.Range("D2").Value = "New York"
$ws.Range("E2").Value = "william Street"
$ws.Range("F2").Value = 10001
$ws.Range("G2").Value = 7013606690

# Set column G width (auto best-fit sized to 11 characters)
$ws.Columns.Item(7).ColumnWidth = 10.14

# Update selection to G2
$ws.Range("G2").Select()
